$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '26.489.41'
$ws.Range("E2").Value = '  +0.90%  '
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '1.726.36'
$ws.Range("E3").Value = '  +0.46%  '
$rng.Style = "Normal"

$rng = $ws.Range("D4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = '244.75'
$ws.Range("E5").Value = '  +2.06%  '
$rng.Style = "Normal"

$rng = $ws.Range("E6")
$rng.NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '
$rng.Style = "Normal"

$rng = $ws.Range("E7")
$rng.NumberFormat = "@"
$ws.Range("E7").Value = '  +1.98%  '
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '0.2670'
$ws.Range("E8").Value = '  +1.80%  '
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = '0.06187'
$ws.Range("E9").Value = '  -0.18%  '
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '1.735.54'
$ws.Range("E10").Value = '  +0.97%  '
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '0.07191'
$ws.Range("E11").Value = '  +1.63%  '
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = '15.58'
$ws.Range("E12").Value = '  +1.71%  '
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = '0.6114'
$ws.Range("E13").Value = '  +2.28%  '
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '4.526'
$ws.Range("E14").Value = '  +2.24%  '
$rng.Style = "Normal"

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = '77.13'
$ws.Range("E15").Value = '  +1.23%  '
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = '0.9997'
$ws.Range("E16").Value = '  -0.06%  '
$rng.Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = '26.501.78'
$ws.Range("E17").Value = '  +0.88%  '
$rng.Style = "Normal"

$rng = $ws.Range("D18:E18")
$rng.NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  -0.09%  '
$rng.Style = "Normal"

$rng = $ws.Range("D19:E19")
$rng.NumberFormat = "@"
$ws.Range("D19").Value = '0.000006947'
$ws.Range("E19").Value = '  +2.12%  '
$rng.Style = "Normal"

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '11.52'
$ws.Range("E20").Value = '  -0.15%  '
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = '1.957.40'
$ws.Range("E21").Value = '  +1.02%  '
$rng.Style = "Normal"

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = '4.516'
$ws.Range("E22").Value = '  -0.50%  '
$rng.Style = "Normal"

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = '8.790'
$ws.Range("E23").Value = '  +0.66%  '
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '5.245'
$ws.Range("E24").Value = '  -0.83%  '
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '137.02'
$ws.Range("E25").Value = '  +1.96%  '
$rng.Style = "Normal"

$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$ws.Range("E26").Value = '  +1.03%  '
$rng.Style = "Normal"

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = '1.777'
$ws.Range("E27").Value = '  +0.79%  '
$rng.Style = "Normal"

$rng = $ws.Range("D28:E28")
$rng.NumberFormat = "@"
$ws.Range("D28").Value = '1.394'
$ws.Range("E28").Value = '  -0.57%  '
$rng.Style = "Normal"

$rng = $ws.Range("D29:E29")
$rng.NumberFormat = "@"
$ws.Range("D29").Value = '107.24'
$ws.Range("E29").Value = '  +0.17%  '
$rng.Style = "Normal"

$rng = $ws.Range("D30:E30")
$rng.NumberFormat = "@"
$ws.Range("D30").Value = '3.962'
$ws.Range("E30").Value = '  -0.31%  '
$rng.Style = "Normal"

$rng = $ws.Range("D31:E31")
$rng.NumberFormat = "@"
$ws.Range("D31").Value = '0.08011'
$ws.Range("E31").Value = '  +3.20%  '
$rng.Style = "Normal"

$rng = $ws.Range("D32:E32")
$rng.NumberFormat = "@"
$ws.Range("D32").Value = '3.690'
$ws.Range("E32").Value = '  +0.42%  '
$rng.Style = "Normal"

$rng = $ws.Range("D33:E33")
$rng.NumberFormat = "@"
$ws.Range("D33").Value = '0.04518'
$ws.Range("E33").Value = '  +1.27%  '
$rng.Style = "Normal"

$rng = $ws.Range("B34:E34")
$rng.NumberFormat = "@"
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.612'
$ws.Range("E34").Value = '  -0.19%  '
$rng.Style = "Normal"

$rng = $ws.Range("B35:E35")
$rng.NumberFormat = "@"
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '0.9969'
$ws.Range("E35").Value = '  +2.16%  '
$rng.Style = "Normal"

$rng = $ws.Range("B36:E36")
$rng.NumberFormat = "@"
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6270'
$ws.Range("E36").Value = '  +1.48%  '
$rng.Style = "Normal"

$rng = $ws.Range("B37:E37")
$rng.NumberFormat = "@"
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9136'
$ws.Range("E37").Value = '  -1.46%  '
$rng.Style = "Normal"

$rng = $ws.Range("B38:E38")
$rng.NumberFormat = "@"
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.079'
$ws.Range("E38").Value = '  +8.26%  '
$rng.Style = "Normal"

$rng = $ws.Range("B39:E39")
$rng.NumberFormat = "@"
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.369'
$ws.Range("E39").Value = '  -2.28%  '
$rng.Style = "Normal"

$rng = $ws.Range("B40:E40")
$rng.NumberFormat = "@"
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '1.001'
$ws.Range("E40").Value = '  +0.11%  '
$rng.Style = "Normal"

$rng = $ws.Range("B41:E41")
$rng.NumberFormat = "@"
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").Value = '103.36'
$ws.Range("E41").Value = '  -9.01%  '
$rng.Style = "Normal"

$rng = $ws.Range("B42:E42")
$rng.NumberFormat = "@"
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.01508'
$ws.Range("E42").Value = '  +1.86%  '
$rng.Style = "Normal"

$rng = $ws.Range("B43:E43")
$rng.NumberFormat = "@"
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.643'
$ws.Range("E43").Value = '  +2.80%  '
$rng.Style = "Normal"

$rng = $ws.Range("B44:E44")
$rng.NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3864'
$ws.Range("E44").Value = '  +1.05%  '
$rng.Style = "Normal"

$rng = $ws.Range("B45:E45")
$rng.NumberFormat = "@"
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.977'
$ws.Range("E45").Value = '  +11.02%  '
$rng.Style = "Normal"

$rng = $ws.Range("B46:E46")
$rng.NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1181'
$ws.Range("E46").Value = '  +0.12%  '
$rng.Style = "Normal"

$rng = $ws.Range("B47:E47")
$rng.NumberFormat = "@"
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05353'
$ws.Range("E47").Value = '  +1.54%  '
$rng.Style = "Normal"

$rng = $ws.Range("B48:E48")
$rng.NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '7.831'
$ws.Range("E48").Value = '  +0.64%  '
$rng.Style = "Normal"

$rng = $ws.Range("D49:E49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = '30.46'
$ws.Range("E49").Value = '  +0.69%  '
$rng.Style = "Normal"

$rng = $ws.Range("B50:E50")
$rng.NumberFormat = "@"
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.251'
$ws.Range("E50").Value = '  +2.90%  '
$rng.Style = "Normal"

$rng = $ws.Range("B51:E51")
$rng.NumberFormat = "@"
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '51.15'
$ws.Range("E51").Value = '  +0.91%  '
$rng.Style = "Normal"
